$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.541.61"
$ws.Range("E2").Value = "  +4.05%  "

$ws.Range("D3").Value = "3.488.99"
$ws.Range("E3").Value = "  +2.72%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.24"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "3.488.59"
$ws.Range("E8").Value = "  +2.77%  "

$ws.Range("E9").Value = "  +7.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.32"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.51%  "

$ws.Range("E11").Value = "  +6.17%  "

$ws.Range("E12").Value = "  +3.46%  "

$ws.Range("D13").Value = "4.092.98"
$ws.Range("E13").Value = "  +2.73%  "

$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.07"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.86%  "

$ws.Range("D16").Value = "66.576.31"
$ws.Range("E16").Value = "  +4.06%  "

$ws.Range("E17").Value = "  +3.05%  "

$ws.Range("D18").Value = "3.500.29"
$ws.Range("E18").Value = "  +3.12%  "

$ws.Range("E19").Value = "  +2.99%  "

$ws.Range("E20").Value = "  +4.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "390.23"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.72%  "

$ws.Range("E22").Value = "  +1.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.81"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.531"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.05%  "

$ws.Range("E26").Value = "  +6.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.29"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +9.36%  "

$ws.Range("E28").Value = "  +1.48%  "

$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("E30").Value = "  +4.35%  "

$ws.Range("E31").Value = "  +4.81%  "

$ws.Range("E32").Value = "  +2.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.53"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.36"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.97%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("E36").Value = "  +8.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.58"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.60%  "

$ws.Range("E38").Value = "  +4.41%  "

$ws.Range("E39").Value = "  +5.70%  "

$ws.Range("E40").Value = "  +5.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0742"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.31"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.68%  "

$ws.Range("D44").Value = "2.784.50"
$ws.Range("E44").Value = "  +2.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.07"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.06%  "

$ws.Range("E46").Value = "  +3.43%  "

$ws.Range("E47").Value = "  +2.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.50"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "345.53"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.74%  "

$ws.Range("E50").Value = "  +4.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.69"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +12.43%  "
